# Weekly price-sheet update ("Fruta / hortaliza, semanal"):
# a new weekly record is inserted as row 27 (most recent date first),
# pushing every existing record down by one row (old row 27 -> new row
# 28, ..., old row 103 -> new row 104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 27:103 down to 28:104, leaving row 27 empty for the new record.
$ws.Rows("27:27").Insert()

# Populate the new record in row 27.
$ws.Range("A27").Value = 8
$ws.Range("B27").Value = "Terminal La Palmera de La Serena"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44764
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 100112052
$ws.Range("G27").Value = "Albahaca"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 4500
$ws.Range("M27").Value = 4250
$ws.Range("N27").Value = "$/paquete"
$ws.Range("O27").Value = "Provincia del Elquí"
$ws.Range("P27").Value = 4250
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"
